$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "OP871U"
$ws.Range("B17").Value = "Cabezal Epson"
$ws.Range("C17").Value = "TMU 220"
$ws.Range("D17").Value = 100000
$ws.Range("E17").Value = 300000
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 14
$ws.Range("H17").Formula = "=(E17-D17)*G17"
$ws.Range("I17").Formula = "=D17*F17"
$ws.Range("J17").Value = 0
